$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.525.75"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.851.13"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").Value = "'241.82"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.6285"
$ws.Range("E6").Value = "  -2.07%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8: OKB
$ws.Range("D8").Value = "'47.87"
$ws.Range("E8").Value = "  +0.98%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.07570"
$ws.Range("E9").Value = "  +0.63%  "

# Row 10: Cardano
$ws.Range("D10").Value = "'0.2976"
$ws.Range("E10").Value = "  -0.12%  "

# Row 11: Solana
$ws.Range("D11").Value = "'24.33"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +0.03%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.881.33"
$ws.Range("E13").Value = "  +1.12%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'5.016"
$ws.Range("E14").Value = "  -0.46%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.6849"
$ws.Range("E15").Value = "  -0.98%  "

# Row 16: Litecoin
$ws.Range("E16").Value = "  -0.02%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "'0.000009808"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18: WrappedliquidstakedEther2.0
$ws.Range("D18").Value = "2.139.01"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19: Uniswap
$ws.Range("E19").Value = "  +1.98%  "

# Row 20: WrappedBTC
$ws.Range("D20").Value = "29.581.47"
$ws.Range("E20").Value = "  -0.60%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'234.72"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22: Avalanche
$ws.Range("E22").Value = "  -1.47%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "'7.620"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25: BinanceUSD
$ws.Range("D25").Value = "'1.000"
$ws.Range("E25").Value = "  +0.04%  "

# Row 26: Monero
$ws.Range("D26").Value = "'155.85"
$ws.Range("E26").Value = "  -1.94%  "

# Row 27: Stellar
$ws.Range("D27").Value = "'0.1390"
$ws.Range("E27").Value = "  -2.12%  "

# Row 28: Cosmos
$ws.Range("D28").Value = "'8.431"
$ws.Range("E28").Value = "  -1.27%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'17.72"

# Row 30: PancakeSwap
$ws.Range("D30").Value = "'1.481"
$ws.Range("E30").Value = "  -0.94%  "

# Row 31: Hedera
$ws.Range("D31").Value = "'0.05836"
$ws.Range("E31").Value = "  -6.43%  "

# Row 32: Toncoin
$ws.Range("D32").Value = "'1.260"
$ws.Range("E32").Value = "  -1.94%  "

# Row 33: Filecoin
$ws.Range("D33").Value = "'4.110"
$ws.Range("E33").Value = "  -1.27%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = "'4.041"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35: LidoDAOToken
$ws.Range("D35").Value = "'1.898"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36: ARBITRUM
$ws.Range("E36").Value = "  +0.03%  "

# Row 37: ImmutableX
$ws.Range("D37").Value = "'0.7171"
$ws.Range("E37").Value = "  -1.55%  "

# Row 38: HuobiToken
$ws.Range("E38").Value = "  -0.59%  "

# Row 39: MXToken
$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "  -0.94%  "

# Row 40: Maker
$ws.Range("D40").Value = "1.235.06"
$ws.Range("E40").Value = "  +2.71%  "

# Row 41: VeChain
$ws.Range("D41").Value = "'0.01778"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42: TrustWalletToken
$ws.Range("D42").Value = "'0.9125"
$ws.Range("E42").Value = "  -1.40%  "

# Row 43: FraxShare
$ws.Range("D43").Value = "'6.137"
$ws.Range("E43").Value = "  -1.68%  "

# Row 44: RocketPoolETH
$ws.Range("D44").Value = "2.049.21"
$ws.Range("E44").Value = "  +0.70%  "

# Row 45: PaxDollar
$ws.Range("D45").Value = "'0.9992"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46: Quant
$ws.Range("D46").Value = "'101.88"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47: Aave
$ws.Range("D47").Value = "'67.46"
$ws.Range("E47").Value = "  +1.43%  "

# Row 48: Aptos
$ws.Range("D48").Value = "'7.287"
$ws.Range("E48").Value = "  +9.02%  "

# Row 49: BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000117"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50: EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.138"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51: TheSandbox
$ws.Range("D51").Value = "'0.4030"
$ws.Range("E51").Value = "  -0.77%  "
